$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_8_5_1"
$ws.Range("B2").Value = -0.3429419163213907
$ws.Range("C2").Value = -0.5818750137294237
$ws.Range("D2").Value = -3.408068147388144
$ws.Range("E2").Value = -0.8700733678094812
$ws.Range("F2").Value = 1.486240983009338
$ws.Range("G2").Value = 2.753446817398071
$ws.Range("H2").Value = 2.644283056259155
$ws.Range("I2").Value = 2.702074766159058

$ws.Range("A3").Value = "model_8_5_0"
$ws.Range("B3").Value = -0.2625370420336641
$ws.Range("C3").Value = -0.5327941198249013
$ws.Range("D3").Value = 0.5820204672620912
$ws.Range("E3").Value = -0.05921902786580735
$ws.Range("F3").Value = 1.39725649356842
$ws.Range("G3").Value = 2.668015480041504
$ws.Range("H3").Value = 0.250734806060791
$ws.Range("I3").Value = 1.530468940734863

$ws.Range("A4").Value = "model_8_5_23"
$ws.Range("B4").Value = 0.006374449395167958
$ws.Range("C4").Value = -0.2683955663146094
$ws.Range("D4").Value = -7.66661357337847
$ws.Range("E4").Value = -1.502148818603772
$ws.Range("F4").Value = 1.099650621414185
$ws.Range("G4").Value = 2.207797288894653
$ws.Range("H4").Value = 5.198871612548828
$ws.Range("I4").Value = 3.615362644195557

$ws.Range("A5").Value = "model_8_5_24"
$ws.Range("B5").Value = 0.006374449395167958
$ws.Range("C5").Value = -0.2683955663146094
$ws.Range("D5").Value = -7.66661357337847
$ws.Range("E5").Value = -1.502148818603772
$ws.Range("F5").Value = 1.099650621414185
$ws.Range("G5").Value = 2.207797288894653
$ws.Range("H5").Value = 5.198871612548828
$ws.Range("I5").Value = 3.615362644195557

$ws.Range("A6").Value = "model_8_5_22"
$ws.Range("B6").Value = 0.006446057720937981
$ws.Range("C6").Value = -0.268227880060512
$ws.Range("D6").Value = -7.66623869028111
$ws.Range("E6").Value = -1.501967933355769
$ws.Range("F6").Value = 1.099571466445923
$ws.Range("G6").Value = 2.207505464553833
$ws.Range("H6").Value = 5.198646545410156
$ws.Range("I6").Value = 3.615101099014282

$ws.Range("A7").Value = "model_8_5_21"
$ws.Range("B7").Value = 0.006600025086089478
$ws.Range("C7").Value = -0.2677690778935686
$ws.Range("D7").Value = -7.665752528365839
$ws.Range("E7").Value = -1.501579700366195
$ws.Range("F7").Value = 1.099401116371155
$ws.Range("G7").Value = 2.206707000732422
$ws.Range("H7").Value = 5.198354721069336
$ws.Range("I7").Value = 3.614540338516235

$ws.Range("A8").Value = "model_8_5_20"
$ws.Range("B8").Value = 0.006616558521031157
$ws.Range("C8").Value = -0.2677182285714299
$ws.Range("D8").Value = -7.665702819788175
$ws.Range("E8").Value = -1.5015380824206
$ws.Range("F8").Value = 1.099382758140564
$ws.Range("G8").Value = 2.206618309020996
$ws.Range("H8").Value = 5.198325157165527
$ws.Range("I8").Value = 3.614480018615723

$ws.Range("A9").Value = "model_8_5_19"
$ws.Range("B9").Value = 0.008583161262560646
$ws.Range("C9").Value = -0.266287818233176
$ws.Range("D9").Value = -7.644961749514746
$ws.Range("E9").Value = -1.496574112965016
$ws.Range("F9").Value = 1.097206234931946
$ws.Range("G9").Value = 2.204128742218018
$ws.Range("H9").Value = 5.185883045196533
$ws.Range("I9").Value = 3.60730767250061

$ws.Range("A10").Value = "model_8_5_16"
$ws.Range("B10").Value = 0.00862164769252538
$ws.Range("C10").Value = -0.2838330103810749
$ws.Range("D10").Value = -7.587054161406236
$ws.Range("E10").Value = -1.496449929488521
$ws.Range("F10").Value = 1.097163796424866
$ws.Range("G10").Value = 2.234668016433716
$ws.Range("H10").Value = 5.151145458221436
$ws.Range("I10").Value = 3.607128381729126

$ws.Range("A11").Value = "model_8_5_18"
$ws.Range("B11").Value = 0.009236903337118019
$ws.Range("C11").Value = -0.2684865965809895
$ws.Range("D11").Value = -7.629292529980599
$ws.Range("E11").Value = -1.494915615077756
$ws.Range("F11").Value = 1.096482872962952
$ws.Range("G11").Value = 2.207955837249756
$ws.Range("H11").Value = 5.176483631134033
$ws.Range("I11").Value = 3.604911327362061

$ws.Range("A12").Value = "model_8_5_17"
$ws.Range("B12").Value = 0.009852864636893499
$ws.Range("C12").Value = -0.2707457837493741
$ws.Range("D12").Value = -7.613963569361731
$ws.Range("E12").Value = -1.493360910248218
$ws.Range("F12").Value = 1.095801115036011
$ws.Range("G12").Value = 2.211888313293457
$ws.Range("H12").Value = 5.167287826538086
$ws.Range("I12").Value = 3.602664947509766

$ws.Range("A13").Value = "model_8_5_15"
$ws.Range("B13").Value = 0.01580042088889533
$ws.Range("C13").Value = -0.2841630586953812
$ws.Range("D13").Value = -7.493086246831236
$ws.Range("E13").Value = -1.478300673459445
$ws.Range("F13").Value = 1.089218854904175
$ws.Range("G13").Value = 2.235242605209351
$ws.Range("H13").Value = 5.094776630401611
$ws.Range("I13").Value = 3.580904006958008

$ws.Range("A14").Value = "model_8_5_13"
$ws.Range("B14").Value = 0.01799467468403748
$ws.Range("C14").Value = -0.3455646204133853
$ws.Range("D14").Value = -7.263025876403258
$ws.Range("E14").Value = -1.472513732045649
$ws.Range("F14").Value = 1.086790442466736
$ws.Range("G14").Value = 2.342119693756104
$ws.Range("H14").Value = 4.956769943237305
$ws.Range("I14").Value = 3.572542428970337

$ws.Range("A15").Value = "model_8_5_14"
$ws.Range("B15").Value = 0.01889947362167255
$ws.Range("C15").Value = -0.3055622453746418
$ws.Range("D15").Value = -7.382423476234868
$ws.Range("E15").Value = -1.470328759726607
$ws.Range("F15").Value = 1.085789084434509
$ws.Range("G15").Value = 2.272490501403809
$ws.Range("H15").Value = 5.028393268585205
$ws.Range("I15").Value = 3.569385766983032

$ws.Range("A16").Value = "model_8_5_12"
$ws.Range("B16").Value = 0.02493503715897361
$ws.Range("C16").Value = -0.4248833315908391
$ws.Range("D16").Value = -6.906772082688796
$ws.Range("E16").Value = -1.453498088030779
$ws.Range("F16").Value = 1.0791095495224
$ws.Range("G16").Value = 2.480183601379395
$ws.Range("H16").Value = 4.743062496185303
$ws.Range("I16").Value = 3.545066833496094

$ws.Range("A17").Value = "model_8_5_7"
$ws.Range("B17").Value = 0.07192435033853706
$ws.Range("C17").Value = 0.22007557707923
$ws.Range("D17").Value = -7.148915132649167
$ws.Range("E17").Value = -1.089474846131153
$ws.Range("F17").Value = 1.027106285095215
$ws.Range("G17").Value = 1.357553839683533
$ws.Range("H17").Value = 4.888318061828613
$ws.Range("I17").Value = 3.019088506698608

$ws.Range("A18").Value = "model_8_5_8"
$ws.Range("B18").Value = 0.08870542090302036
$ws.Range("C18").Value = -1.273536409629591
$ws.Range("D18").Value = -2.906810198756735
$ws.Range("E18").Value = -1.213258733535874
$ws.Range("F18").Value = 1.008534669876099
$ws.Range("G18").Value = 3.957367897033691
$ws.Range("H18").Value = 2.343591690063477
$ws.Range("I18").Value = 3.197944164276123

$ws.Range("A19").Value = "model_8_5_11"
$ws.Range("B19").Value = 0.1290143972237719
$ws.Range("C19").Value = -0.04663699373921926
$ws.Range("D19").Value = -6.675924926984043
$ws.Range("E19").Value = -1.167166063312758
$ws.Range("F19").Value = 0.9639244079589844
$ws.Range("G19").Value = 1.821799635887146
$ws.Range("H19").Value = 4.604583740234375
$ws.Range("I19").Value = 3.13134503364563

$ws.Range("A20").Value = "model_8_5_6"
$ws.Range("B20").Value = 0.1497690007007507
$ws.Range("C20").Value = 0.2483970604219117
$ws.Range("D20").Value = -6.01112311479076
$ws.Range("E20").Value = -0.8491210798624353
$ws.Range("F20").Value = 0.9409551620483398
$ws.Range("G20").Value = 1.30825662612915
$ws.Range("H20").Value = 4.20578670501709
$ws.Range("I20").Value = 2.67180061340332

$ws.Range("A21").Value = "model_8_5_2"
$ws.Range("B21").Value = 0.1567435037229445
$ws.Range("C21").Value = 0.151144403821832
$ws.Range("D21").Value = -3.220768026440485
$ws.Range("E21").Value = -0.3659874385194644
$ws.Range("F21").Value = 0.9332365393638611
$ws.Range("G21").Value = 1.477536797523499
$ws.Range("H21").Value = 2.53192663192749
$ws.Range("I21").Value = 1.973719477653503

$ws.Range("A22").Value = "model_8_5_10"
$ws.Range("B22").Value = 0.1666937932430046
$ws.Range("C22").Value = -0.4665719639325643
$ws.Range("D22").Value = -4.801711666323968
$ws.Range("E22").Value = -1.068816758620074
$ws.Range("F22").Value = 0.9222244620323181
$ws.Range("G22").Value = 2.552747488021851
$ws.Range("H22").Value = 3.480292797088623
$ws.Range("I22").Value = 2.989239692687988

$ws.Range("A23").Value = "model_8_5_5"
$ws.Range("B23").Value = 0.2129180315444903
$ws.Range("C23").Value = 0.2721359770476656
$ws.Range("D23").Value = -5.022828561135464
$ws.Range("E23").Value = -0.6408945933471351
$ws.Range("F23").Value = 0.8710677027702332
$ws.Range("G23").Value = 1.266936302185059
$ws.Range("H23").Value = 3.612935066223145
$ws.Range("I23").Value = 2.370933532714844

$ws.Range("A24").Value = "model_8_5_9"
$ws.Range("B24").Value = 0.2163355637421781
$ws.Range("C24").Value = -0.26922186295789
$ws.Range("D24").Value = -4.783400621162458
$ws.Range("E24").Value = -0.9393756321024347
$ws.Range("F24").Value = 0.8672856092453003
$ws.Range("G24").Value = 2.209235668182373
$ws.Range("H24").Value = 3.469308376312256
$ws.Range("I24").Value = 2.802209854125977

$ws.Range("A25").Value = "model_8_5_4"
$ws.Range("B25").Value = 0.2254744745163408
$ws.Range("C25").Value = 0.2846562770626626
$ws.Range("D25").Value = -4.511488917473241
$ws.Range("E25").Value = -0.53300824341172
$ws.Range("F25").Value = 0.8571715354919434
$ws.Range("G25").Value = 1.245143055915833
$ws.Range("H25").Value = 3.306195735931396
$ws.Range("I25").Value = 2.215048313140869

$ws.Range("A26").Value = "model_8_5_3"
$ws.Range("B26").Value = 0.2306290605062832
$ws.Range("C26").Value = 0.1017115171228792
$ws.Range("D26").Value = -3.678599795336724
$ws.Range("E26").Value = -0.4869612302651669
$ws.Range("F26").Value = 0.8514668345451355
$ws.Range("G26").Value = 1.563580870628357
$ws.Range("H26").Value = 2.806567668914795
$ws.Range("I26").Value = 2.148514986038208
